$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278, pushing existing rows 278:296 down to 279:297
$ws.Rows.Item(278).Insert()

# Copy the date number format (style) from the row below (now row 279, originally row 278)
$ws.Range("D278").NumberFormat = $ws.Range("D279").NumberFormat

# Fill in the new row's data
$ws.Cells.Item(278, 1).Value = 4
$ws.Cells.Item(278, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value = "Los Lagos"
$ws.Cells.Item(278, 4).Value = 45021
$ws.Cells.Item(278, 5).Value = 10
$ws.Cells.Item(278, 6).Value = 100112009
$ws.Cells.Item(278, 7).Value = "Acelga"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 15
$ws.Cells.Item(278, 11).Value = 9000
$ws.Cells.Item(278, 12).Value = 9000
$ws.Cells.Item(278, 13).Value = 9000
$ws.Cells.Item(278, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(278, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(278, 16).Value = 750
$ws.Cells.Item(278, 17).Value = 12
$ws.Cells.Item(278, 18).Value = "Hortaliza"
